$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3.283333333333331
$ws.Range("F2").Value = 3
$ws.Range("H2").Value = 2

# Row 23
$ws.Range("D23").Value = 0
$ws.Range("J23").Value = "Out"

# Row 29
$ws.Range("H29").Value = 5

# Row 69
$ws.Range("D69").Value = 16.68333333333333
$ws.Range("J69").ClearContents()

# Row 72
$ws.Range("D72").Value = 14.13333333333333
$ws.Range("F72").Value = 14.25
$ws.Range("H72").Value = 1

# Row 83
$ws.Range("D83").Value = 3.683333333333333
$ws.Range("F83").Value = 3.8
$ws.Range("H83").Value = 2

# Row 102
$ws.Range("D102").Value = -2.216666666666669
$ws.Range("J102").ClearContents()

# Row 104
$ws.Range("H104").Value = 4

# Row 116
$ws.Range("H116").Value = 6

# Row 123
$ws.Range("D123").Value = 0
$ws.Range("J123").Value = "Out"

# Row 137
$ws.Range("D137").Value = 0
$ws.Range("J137").Value = "Out"

# Row 151
$ws.Range("D151").Value = 0
$ws.Range("J151").Value = "Out"

# Row 154
$ws.Range("D154").Value = 0
$ws.Range("J154").Value = "Out"

# Row 165
$ws.Range("H165").Value = 5

# Row 166
$ws.Range("H166").Value = 3

# Row 180
$ws.Range("D180").Value = 7.528333333333333
$ws.Range("F180").Value = 7.65
$ws.Range("H180").Value = 1

# Row 186
$ws.Range("H186").Value = 1

# Row 192
$ws.Range("H192").Value = 6

# Row 201
$ws.Range("D201").Value = 14.58333333333333
$ws.Range("F201").Value = 14.5
$ws.Range("H201").Value = 1

# Row 202
$ws.Range("D202").Value = 7.933333333333332
$ws.Range("J202").ClearContents()

# Row 235
$ws.Range("D235").Value = 13.33333333333333
$ws.Range("F235").Value = 13.05
$ws.Range("H235").Value = 1

# Row 243
$ws.Range("H243").Value = 5

# Row 279
$ws.Range("D279").Value = 3.533333333333331
$ws.Range("J279").ClearContents()

# Row 286
$ws.Range("H286").Value = 6

# Row 294
$ws.Range("H294").Value = 5

# Row 295
$ws.Range("D295").Value = 7.383333333333333
$ws.Range("J295").ClearContents()

# Row 299
$ws.Range("H299").Value = 1

# Row 314
$ws.Range("H314").Value = 4

# Row 322
$ws.Range("D322").Value = 8.383333333333333
$ws.Range("F322").Value = 8.300000000000001
$ws.Range("H322").Value = 3

# Row 351
$ws.Range("D351").Value = 14.38333333333333
$ws.Range("J351").ClearContents()

# Row 352
$ws.Range("H352").Value = 1

# Row 402
$ws.Range("D402").Value = 14.73333333333333
$ws.Range("J402").ClearContents()

# Row 410
$ws.Range("F410").Value = 9.550000000000001
$ws.Range("H410").Value = 4

# Row 418
$ws.Range("D418").Value = 7.198333333333331
$ws.Range("J418").Value = "Questionable"

# Row 436
$ws.Range("H436").Value = 5

# Row 457
$ws.Range("H457").Value = 3
